$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Unprotect()
$ws1.Range("C2").Value = 43710
$ws1.Range("D2").Value = 43724
$ws1.Range("E2").Value = 43738
$ws1.Range("F2").Value = 43752
$ws1.Range("G2").Value = 43766
$ws1.Range("H2").Value = 43780
$ws1.Range("I2").Value = 43794
$ws1.Range("J2").Value = 43808
$ws1.Range("C3").ClearContents()
$ws1.Range("J3").Select()
